# Generate Report for Handoff
# Appends a new tracked file (c764132e-79ae-43f4-8896-7956347324a1.md, status
# "Ready for handoff") as row 9 on the Overview / zh-cn / de-de sheets, and
# grows each sheet's table/autofilter/dimension to match.

$wb = $excel.ActiveWorkbook

$fileId   = "c764132e-79ae-43f4-8896-7956347324a1"
$fileName = "$fileId.md"
$filePath = "e2e\$fileId.md"
$ghHash   = "93a866ebaf7171031a99982c782984b787c4d52c"
$ghBase   = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$ghHash"

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Overview sheet: new row 9
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A9").Value = $fileName
$wsOverview.Range("B9").Value = $filePath
$wsOverview.Range("C9").Value = ".md"
$wsOverview.Range("D9").Value = ""
$wsOverview.Range("E9").Value = "Ready for handoff"
$wsOverview.Range("F9").Value = "Ready for handoff"
$wsOverview.Range("G9").NumberFormat = $dateFmt
$wsOverview.Range("G9").Value = "2016-08-28 00:41:46"

$null = $wsOverview.Hyperlinks.Add($wsOverview.Range("B9"), "$ghBase/e2e/$fileName", "", "", $filePath)
$wsOverview.Range("B9").Style = "Hyperlink"

$tblOverview = $wsOverview.ListObjects.Item(1)
$tblOverview.Resize($wsOverview.Range("A1:G9"))

# ---------------------------------------------------------------------------
# zh-cn sheet: new row 9
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A9").Value = $fileName
$wsZhCn.Range("B9").Value = ".md"
$wsZhCn.Range("C9").Value = "Ready for handoff"
$wsZhCn.Range("D9").Value = "e2e"
$wsZhCn.Range("E9").Value = "ht"
$wsZhCn.Range("F9").Value = "False"
$wsZhCn.Range("G9").Value = "$fileId.$ghHash.zh-cn.xlf"
$wsZhCn.Range("H9").NumberFormat = $dateFmt
$wsZhCn.Range("H9").Value = "2016-08-28 00:41:42"
$wsZhCn.Range("I9").Value = ""
$wsZhCn.Range("J9").Value = ""
$wsZhCn.Range("K9").NumberFormat = $dateFmt
$wsZhCn.Range("K9").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("L9").Value = ""
$wsZhCn.Range("M9").Value = "True"
$wsZhCn.Range("N9").Value = ""
$wsZhCn.Range("O9").Value = "False"
$wsZhCn.Range("P9").Value = ""

$null = $wsZhCn.Hyperlinks.Add($wsZhCn.Range("A9"), "$ghBase/e2e/$fileName", "", "", $fileName)
$wsZhCn.Range("A9").Style = "Hyperlink"

$tblZhCn = $wsZhCn.ListObjects.Item(1)
$tblZhCn.Resize($wsZhCn.Range("A1:P9"))

# ---------------------------------------------------------------------------
# de-de sheet: new row 9
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A9").Value = $fileName
$wsDeDe.Range("B9").Value = ".md"
$wsDeDe.Range("C9").Value = "Ready for handoff"
$wsDeDe.Range("D9").Value = "e2e"
$wsDeDe.Range("E9").Value = "ht"
$wsDeDe.Range("F9").Value = "False"
$wsDeDe.Range("G9").Value = "$fileId.$ghHash.de-de.xlf"
$wsDeDe.Range("H9").NumberFormat = $dateFmt
$wsDeDe.Range("H9").Value = "2016-08-28 00:41:46"
$wsDeDe.Range("I9").Value = ""
$wsDeDe.Range("J9").Value = ""
$wsDeDe.Range("K9").NumberFormat = $dateFmt
$wsDeDe.Range("K9").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("L9").Value = ""
$wsDeDe.Range("M9").Value = "True"
$wsDeDe.Range("N9").Value = ""
$wsDeDe.Range("O9").Value = "False"
$wsDeDe.Range("P9").Value = ""

$null = $wsDeDe.Hyperlinks.Add($wsDeDe.Range("A9"), "$ghBase/e2e/$fileName", "", "", $fileName)
$wsDeDe.Range("A9").Style = "Hyperlink"

$tblDeDe = $wsDeDe.ListObjects.Item(1)
$tblDeDe.Resize($wsDeDe.Range("A1:P9"))

Write-Output "Handoff report row added for $fileName"
